# daily auto push: 2025-10-02 06:39 UTC
# Append the day's new data row (2025/10/02, 木, 14, 201) to the bottom of
# the log table on the active sheet (row 51, right after the existing
# row 50), extending the used range from A1:D50 to A1:D51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 51

# Column A holds date-looking text (e.g. "2025/09/22") that must stay a
# literal string rather than being auto-converted to a date serial by
# Excel's input parsing. Temporarily mark the cell as Text before writing
# the value, then restore the default "Normal" style so the new row ends
# up unstyled -- matching the rest of the data rows.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/02"
$cellA.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "木"
$ws.Cells.Item($newRow, 3).Value = 14
$ws.Cells.Item($newRow, 4).Value = 201
